$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1950026.9
$ws.Range("I33").Value = 2179389.5
$ws.Range("K33").Value = 2179389.5
$ws.Range("M33").Value = -2179160.5

$ws.Range("H34").Value = 6173.8
$ws.Range("I34").Value = 1891.1428
$ws.Range("K34").Value = 1891.1428
$ws.Range("M34").Value = -1688.1428

$ws.Range("H36").Value = 6173.8
$ws.Range("I36").Value = 1891.1428
$ws.Range("K36").Value = 1891.1428
$ws.Range("M36").Value = -1176.1428

$ws.Range("H100").Value = 2657.111
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H116").Value = 5833.7144
$ws.Range("J116").Value = 5973.5
$ws.Range("L116").Value = 5973.5
$ws.Range("N116").Value = -12857.5

$ws.Range("H135").Value = 1616.6897
$ws.Range("I135").Value = 695
$ws.Range("K135").Value = 6255
$ws.Range("M135").Value = -3720

$ws.Range("H137").Value = 60098.484
$ws.Range("I137").Value = 128596.36
$ws.Range("K137").Value = 385789.08
$ws.Range("M137").Value = -383239.08

$ws.Range("H138").Value = 1964.75
$ws.Range("I138").Value = 1102.4667
$ws.Range("J138").Value = 2334.3
$ws.Range("K138").Value = 3307.4001
$ws.Range("L138").Value = 7002.900000000001
$ws.Range("M138").Value = 1832.5999
$ws.Range("N138").Value = -17282.9

$ws.Range("H141").Value = 2677.1428
$ws.Range("I141").Value = 2297.5
$ws.Range("J141").Value = 3183.3333
$ws.Range("K141").Value = 6892.5
$ws.Range("L141").Value = 9549.999899999999
$ws.Range("M141").Value = -1712.5
$ws.Range("N141").Value = -19909.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H45").Value = 9593995
$ws.Range("I45").Value = 15985505
$ws.Range("J45").Value = 6731.3335
$ws.Range("K45").Value = 15985505
$ws.Range("L45").Value = 6731.3335
$ws.Range("M45").Value = -15985128
$ws.Range("N45").Value = -7485.3335

$ws.Range("H110").Value = 27777776
$ws.Range("I110").Value = 27777776
$ws.Range("K110").Value = 27777776
$ws.Range("M110").Value = -27775731

$ws.Range("H122").Value = 673834.5
$ws.Range("I122").Value = 1630.8
$ws.Range("J122").Value = 3474683.2
$ws.Range("K122").Value = 4892.4
$ws.Range("L122").Value = 10424049.6
$ws.Range("M122").Value = -2442.4
$ws.Range("N122").Value = -10428949.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3378460.8
$ws.Range("I94").Value = 5682827.5
$ws.Range("K94").Value = 5682827.5
$ws.Range("M94").Value = -5682376.5

$ws.Range("H107").Value = 3760552.8
$ws.Range("I107").Value = 4465429
$ws.Range("K107").Value = 4465429
$ws.Range("M107").Value = -4463509

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3600.3333
$ws.Range("I6").Value = 3900.5
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 3900.5
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -3787.5
$ws.Range("N6").Value = -3226

$ws.Range("H31").Value = 39019.96
$ws.Range("I31").Value = 1786.8
$ws.Range("J31").Value = 63842.066
$ws.Range("K31").Value = 1786.8
$ws.Range("L31").Value = 63842.066
$ws.Range("M31").Value = -1491.8
$ws.Range("N31").Value = -64432.066

$ws.Range("H34").Value = 39019.96
$ws.Range("I34").Value = 1786.8
$ws.Range("J34").Value = 63842.066
$ws.Range("K34").Value = 1786.8
$ws.Range("L34").Value = 63842.066
$ws.Range("M34").Value = -1584.8
$ws.Range("N34").Value = -64246.066

$ws.Range("H58").Value = 2050.6829
$ws.Range("I58").Value = 1769.0322
$ws.Range("J58").Value = 2923.8
$ws.Range("K58").Value = 1769.0322
$ws.Range("L58").Value = 2923.8
$ws.Range("M58").Value = -1566.0322
$ws.Range("N58").Value = -3329.8

$ws.Range("H99").Value = 3142.9412
$ws.Range("I99").Value = 2429.875
$ws.Range("K99").Value = 2429.875
$ws.Range("M99").Value = -931.875

$ws.Range("H105").Value = 1346.4445
$ws.Range("I105").Value = 1352.875
$ws.Range("K105").Value = 1352.875
$ws.Range("M105").Value = 394.125

$ws.Range("H107").Value = 1472.72
$ws.Range("J107").Value = 1697.1666
$ws.Range("L107").Value = 1697.1666
$ws.Range("N107").Value = -5537.1666

$ws.Range("H122").Value = 2764.875
$ws.Range("I122").Value = 2474.4
$ws.Range("J122").Value = 3249
$ws.Range("K122").Value = 7423.200000000001
$ws.Range("L122").Value = 9747
$ws.Range("M122").Value = -4973.200000000001
$ws.Range("N122").Value = -14647

$ws.Range("H126").Value = 3142.9412
$ws.Range("I126").Value = 2429.875
$ws.Range("K126").Value = 7289.625
$ws.Range("M126").Value = -4819.625

$ws.Range("H132").Value = 25592.297
$ws.Range("I132").Value = 1371.2667
$ws.Range("J132").Value = 129396.71
$ws.Range("K132").Value = 4113.800099999999
$ws.Range("L132").Value = 388190.13
$ws.Range("M132").Value = -1583.800099999999
$ws.Range("N132").Value = -393250.13

$ws.Range("H136").Value = 2050.6829
$ws.Range("I136").Value = 1769.0322
$ws.Range("J136").Value = 2923.8
$ws.Range("K136").Value = 5307.096600000001
$ws.Range("L136").Value = 8771.400000000001
$ws.Range("M136").Value = -2757.096600000001
$ws.Range("N136").Value = -13871.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9264899
$ws.Range("I56").Value = 9264899
$ws.Range("K56").Value = 9264899
$ws.Range("M56").Value = -9264369

$ws.Range("H122").Value = 881.2414
$ws.Range("I122").Value = 865.36365
$ws.Range("K122").Value = 7788.27285
$ws.Range("M122").Value = -5338.27285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1254550.5
$ws.Range("I97").Value = 1401903.5
$ws.Range("K97").Value = 1401903.5
$ws.Range("M97").Value = -1401407.5

$ws.Range("H122").Value = 290167.72
$ws.Range("I122").Value = 425594.16
$ws.Range("J122").Value = 5772.2
$ws.Range("K122").Value = 1276782.48
$ws.Range("L122").Value = 17316.6
$ws.Range("M122").Value = -1274332.48
$ws.Range("N122").Value = -22216.6

$ws.Range("H126").Value = 5856959
$ws.Range("I126").Value = 2843839.2
$ws.Range("K126").Value = 8531517.600000001
$ws.Range("M126").Value = -8529047.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 15000001
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws.Range("H46").Value = 5077.7856
$ws.Range("I46").Value = 3962.0667
$ws.Range("J46").Value = 6365.154
$ws.Range("K46").Value = 3962.0667
$ws.Range("L46").Value = 6365.154
$ws.Range("M46").Value = -3774.0667
$ws.Range("N46").Value = -6741.154

$ws.Range("H122").Value = 4242.391
$ws.Range("I122").Value = 2911
$ws.Range("K122").Value = 8733
$ws.Range("M122").Value = -6283

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1852.0625
$ws.Range("I122").Value = 1425.8889
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 4277.6667
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -1827.6667
$ws.Range("N122").Value = -12100

$ws.Range("H132").Value = 22245922
$ws.Range("I132").Value = 24394254
$ws.Range("J132").Value = 225521.75
$ws.Range("K132").Value = 73182762
$ws.Range("L132").Value = 676565.25
$ws.Range("M132").Value = -73180232
$ws.Range("N132").Value = -681625.25
